# Neue Kategorie Ausgaben: Personalaufwand
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ausgaben")

# Rename the "Lohn" category (shared string already used by A16) to
# "Personalaufwand", and re-point the "Sonstige Ausgaben" rows (17-25) at
# the same new category, matching row 16.
$ws.Range("A16").Value = "Personalaufwand"
$ws.Range("A17:A25").Value = "Personalaufwand"

# Move/restore the active selection on the sheet.
$ws.Range("A22").Select()
